$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = "P27"
$ws.Range("G6").Value = 2
$ws.Range("F7").Value = "P28"
$ws.Range("G7").Value = 2
$ws.Range("F8").Value = "P29"
$ws.Range("G8").Value = 2
$ws.Range("F9").Value = "P30"
$ws.Range("G9").Value = 2
$ws.Range("F10").Value = "P31"
$ws.Range("G10").Value = 2
$ws.Range("F11").Value = "P34"
$ws.Range("G11").Value = 2
$ws.Range("F12").Value = "P35"
$ws.Range("G12").Value = 2
$ws.Range("F13").Value = "P36"
$ws.Range("G13").Value = 2
$ws.Range("F14").Value = "P5"
$ws.Range("G14").HorizontalAlignment = -4131
$ws.Range("G14").Value = 3
$ws.Range("F15").Value = "P73"
$ws.Range("G15").HorizontalAlignment = -4131
$ws.Range("G15").Value = 1
$ws.Range("F16").Value = "P3"
$ws.Range("G16").Value = 3
$ws.Range("F17").Value = "P4"
$ws.Range("G17").Value = 3
$ws.Range("F18").Value = "P89"
$ws.Range("G18").Value = 0
$ws.Range("F19").Value = "P90"
$ws.Range("G19").Value = 0
$ws.Range("F20").Value = "P6"
$ws.Range("G20").Value = 3
$ws.Range("F21").Value = "P64"
$ws.Range("G21").Value = 1
$ws.Range("F22").Value = "P72"
$ws.Range("G22").Value = 1
$ws.Range("F23").Value = "P65"
$ws.Range("G23").Value = 1
$ws.Range("F24").Value = "P57"
$ws.Range("G24").Value = 1
$ws.Range("F25").Value = "P9"
$ws.Range("G25").Value = 3
$ws.Range("F26").Value = "P10"
$ws.Range("G26").Value = 3
$ws.Range("F27").Value = "P12"
$ws.Range("G27").Value = 3
$ws.Range("F28").Value = "P13"
$ws.Range("G28").Value = 3
$ws.Range("F29").Value = "P15"
$ws.Range("G29").Value = 3
$ws.Range("F30").Value = "P16"
$ws.Range("G30").Value = 3
$ws.Range("F31").Value = "P71"
$ws.Range("G31").Value = 1
$ws.Range("F32").Value = "P62"
$ws.Range("G32").Value = 1
$ws.Range("F33").Value = "P61"
$ws.Range("G33").Value = 1
$ws.Range("F34").Value = "P60"
$ws.Range("G34").Value = 1
$ws.Range("F35").Value = "P59"
$ws.Range("G35").Value = 1

$ws.Range("I8").Select()
